$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 45998
$ws.Range("B4").Value = 48994

$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat
$ws.Range("B4").NumberFormat = $ws.Range("B3").NumberFormat
